# Coalesce the runs that make up the Title paragraph ("Testing" + " " +
# "custom" + " " + "properties") and the Author paragraph ("A." + " " +
# "M.") into single runs, matching the docx writer's new behaviour of
# emitting one run per contiguous span of uniformly-formatted text
# instead of a separate <w:r> for every token/space.
#
# We do this the way a human editing in Word would: keep the paragraph's
# first run in place (so its original formatting/whitespace markers are
# preserved), delete the remaining runs' text, then type the remaining
# words back in right after the first run. That merges everything into
# one run per paragraph without disturbing anything else in the
# document.

$d = $word.ActiveDocument

# --- Paragraph 1 (style "Title"): "Testing custom properties" ---
$p1 = $d.Paragraphs.Item(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End - 1   # exclude the paragraph mark

$firstWordLen = 7   # Len("Testing")
$run1 = $d.Range($p1Start, $p1Start + $firstWordLen)
$rest1 = $d.Range($p1Start + $firstWordLen, $p1End)
$rest1.Delete()
$run1.InsertAfter(" custom properties")

# --- Paragraph 2 (style "Author"): "A. M." ---
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1   # exclude the paragraph mark

$firstWordLen2 = 2   # Len("A.")
$run2 = $d.Range($p2Start, $p2Start + $firstWordLen2)
$rest2 = $d.Range($p2Start + $firstWordLen2, $p2End)
$rest2.Delete()
$run2.InsertAfter(" M.")
